# Delete the data row for account 004565146 / GUSTAVO / 25410.86
# which is the 4th row (including the header row) on the "Export" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

$ws.Rows.Item(4).Delete()
